$d = $word.ActiveDocument

# "personel ${pangkat_terlapor} ${nama_terlapor} jabatan" occurs exactly
# once in the whole document (in the "ditemukan cukup bukti ..." /
# dugaan pelanggaran paragraph), so this Find uniquely anchors the
# target run -- unlike a bare "${nama_terlapor} jabatan" which also
# matches the earlier, unrelated "Sehubungan dengan rujukan ..."
# paragraph.
$rng = $d.Content
$found = $rng.Find.Execute('personel ${pangkat_terlapor} ${nama_terlapor} jabatan', $false, $false, $false, $false, $false, $true, 1, $false, "", 0)

# "personel ${pangkat_terlapor} ${nama_terlapor" is 44 characters long,
# so the closing brace + trailing space of the nama_terlapor placeholder
# start right after it.
$braceStart = $rng.Start + 44
$braceRange = $d.Range($braceStart, $braceStart + 2)   # covers "} "

# Split the single "} " run into three runs: "}", " ", " ".
# Step 1: shrink this run down to just the closing brace.
$braceRange.Text = "}"

# Step 2: insert the first space right after it, then nudge formatting on
# the (now brace-only) preceding run so the two stay distinct runs
# instead of being re-merged back together.
$insPt1 = $d.Range($braceRange.End, $braceRange.End)
$insPt1.InsertAfter(" ")
$braceRange.Bold = 1
$braceRange.Bold = 0

# Step 3: insert the second space after the first, nudging the first
# space run the same way so all three remain separate runs.
$insPt2 = $d.Range($insPt1.End, $insPt1.End)
$insPt2.InsertAfter(" ")
$insPt1.Bold = 1
$insPt1.Bold = 0
